{"js": "// Append \" (Changed main)\" after \"This is a Microsoft word document.\" as\n// three separate runs (\" (\", \"Changed main\", \")\") \u2014 matching the target\n// OOXML diff, which shows the original run left untouched and three new\n// sibling <w:r> elements appended to the same paragraph.\n\nconst body = context.document.body;\n\n// Locate the exact sentence so the edit lands correctly regardless of the\n// paragraph's position in the document.\nconst results = body.search(\"This is a Microsoft word document.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the target sentence to edit.\");\n}\n\nconst target = results.items[0];\n\n// Use a flat-OPC OOXML fragment so the three new runs stay distinct\n// (plain insertText() calls get merged into the existing run since they\n// share identical formatting).\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n            <w:r><w:t>Changed main</w:t></w:r>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(flatOpcXml, \"End\");\nawait context.sync();\n", "ps1": "# Append \" (Changed main)\" after \"This is a Microsoft word document.\" as\n# three separate runs (\" (\", \"Changed main\", \")\") \u2014 matching the target\n# OOXML diff, which shows the original run left untouched and three new\n# sibling <w:r> elements appended to the same paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the exact sentence via Find so the edit lands correctly regardless\n# of the paragraph's position in the document. The returned range spans\n# exactly the matched text (it does not include the paragraph mark).\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"This is a Microsoft word document.\")\nif (-not $found) {\n    throw \"Could not find the target sentence to edit.\"\n}\n\n# Use a flat-OPC OOXML fragment so the three new runs stay distinct (plain\n# InsertAfter() text merges into the existing run since it shares identical\n# formatting). Passing \"End\" as the insertion location appends the runs as\n# siblings inside the matched paragraph instead of replacing the range's\n# content or minting a new paragraph.\n$xml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n            <w:r><w:t>Changed main</w:t></w:r>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($xml, \"End\")\n"}
